# Applies the "more datacamp and disk_savvy data" edit:
#   - Inserts a new row (course entry) "Data Analysis in Excel" with a rating
#     of 5 in column K, right after the existing "Introduction to data" row.
#   - Leaves one blank (but styled) spacer row below it, matching the sheet's
#     existing layout convention, then the remaining course table shifts down
#     by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 47 is currently the spacer row between the "skills" mini-table (rows
# 1-46) and the "courses" table (rows 48-53). Insert a fresh row there and
# push everything below it (the whole courses table) down by one.
$ws.Range("A47:K47").Insert(-4121) | Out-Null

# The insert copies the formatting of the row above (row 46) across every
# column that had an explicit style there (I and K); only A and K should
# carry a style on the new row, so drop the stray formatted-but-empty I47.
$ws.Cells.Item(47, 9).Clear()

# Give the new row the same text style (s="3", the dark "entry title" font)
# used by the other standalone entries like A46/A53/A54.
$ws.Range("A46").Copy()
$ws.Range("A47").PasteSpecial(-4122) | Out-Null

# Row 48 becomes an empty spacer row that still carries style 3 on A48 (as
# in the target), matching the pattern used elsewhere in the sheet.
$ws.Range("A47").Copy()
$ws.Range("A48").PasteSpecial(-4122) | Out-Null

# Fill in the new entry's data.
$ws.Cells.Item(47, 1).Value2 = "Data Analysis in Excel"
$ws.Cells.Item(47, 11).Value2 = 5

# Match the final cursor position recorded in the saved workbook.
$ws.Range("A55").Select() | Out-Null

$excel.CutCopyMode = 0
